$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "28.456.08"
$ws.Range("E2").Value = "  +0.73%  "

# Row 3
$ws.Range("D3").Value = "1.915.28"
$ws.Range("E3").Value = "  +2.27%  "

# Row 4
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.22%  "

# Row 6
$ws.Range("E6").Value = "  +0.08%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5109"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.71%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3964"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.64%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09711"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.47%  "

# Row 10
$ws.Range("E10").Value = "  +1.80%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.13"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.97%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.466"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.14%  "

# Row 13
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.933.35"
$ws.Range("E13").Value = "  +3.02%  "

# Row 14
$ws.Range("B14").Value = "Solana"
$ws.Range("C14").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.03"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.58%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.395"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.18%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.002"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.11%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001130"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.55%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "93.82"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.44%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06678"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.25%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.09"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.98%  "

# Row 21
$ws.Range("E21").Value = "  +0.03%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.262"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.76%  "

# Row 23
$ws.Range("D23").Value = "28.514.31"
$ws.Range("E23").Value = "  +0.67%  "

# Row 24
$ws.Range("E24").Value = "  +1.16%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.319"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.60%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.669"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +5.56%  "

# Row 27
$ws.Range("D27").Value = "2.137.62"
$ws.Range("E27").Value = "  +2.38%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "21.18"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.57%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "158.54"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.38%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "128.36"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.51%  "

# Row 31
$ws.Range("E31").Value = "  +4.30%  "

# Row 32
$ws.Range("E32").Value = "  +0.65%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.694"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.02%  "

# Row 34
$ws.Range("E34").Value = "  +1.02%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.830"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.72%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06710"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.46%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02439"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.10%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.253"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.49%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2221"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.69%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "11.63"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.50%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6425"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.17%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.056"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.87%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.212"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.15%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.002"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.08%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.53"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.15%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6061"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.87%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.786"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.28%  "

# Row 48
$ws.Range("E48").Value = "  +0.02%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.053"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.15%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "124.96"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.19%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.196"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.08%  "
